# Updated cryptos list on Mon May 20 09:37:38 UTC 2024 with GitHub Actions
#
# Refreshes the coin Price (D) and Volume(1h) (E) columns with a newer
# snapshot from coinranking.com. A few rows also swapped rank with their
# neighbor, so their Coin name (B) and Link (C) move together with the
# row's new Price/Volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "66.960.98" }
    @{ Cell = "E2"; Value = "  -0.13%  " }
    @{ Cell = "D3"; Value = "3.104.75" }
    @{ Cell = "E3"; Value = "  +0.11%  " }
    @{ Cell = "E4"; Value = "  +0.03%  " }
    @{ Cell = "D5"; Value = "576.41" }
    @{ Cell = "E5"; Value = "  -0.45%  " }
    @{ Cell = "D6"; Value = "177.48" }
    @{ Cell = "E6"; Value = "  +2.99%  " }
    @{ Cell = "E7"; Value = "  -0.06%  " }
    @{ Cell = "D8"; Value = "3.101.79" }
    @{ Cell = "E8"; Value = "  +0.15%  " }
    @{ Cell = "E9"; Value = "  -1.03%  " }
    @{ Cell = "E10"; Value = "  -2.18%  " }
    @{ Cell = "D11"; Value = "0.152" }
    @{ Cell = "E11"; Value = "  -0.41%  " }
    @{ Cell = "D12"; Value = "0.467" }
    @{ Cell = "E12"; Value = "  -1.72%  " }
    @{ Cell = "E13"; Value = "  -2.71%  " }
    @{ Cell = "E14"; Value = "  -1.31%  " }
    @{ Cell = "D16"; Value = "3.624.42" }
    @{ Cell = "E16"; Value = "  +0.27%  " }
    @{ Cell = "D17"; Value = "66.926.02" }
    @{ Cell = "E17"; Value = "  -0.10%  " }
    @{ Cell = "D18"; Value = "7.03" }
    @{ Cell = "E18"; Value = "  -0.53%  " }
    @{ Cell = "B19"; Value = "WrappedEther" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D19"; Value = "3.104.71" }
    @{ Cell = "E19"; Value = "  +0.08%  " }
    @{ Cell = "B20"; Value = "Chainlink" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" }
    @{ Cell = "D20"; Value = "16.74" }
    @{ Cell = "E20"; Value = "  +1.49%  " }
    @{ Cell = "D21"; Value = "479.54" }
    @{ Cell = "E21"; Value = "  -1.85%  " }
    @{ Cell = "D22"; Value = "7.77" }
    @{ Cell = "E22"; Value = "  -0.37%  " }
    @{ Cell = "D23"; Value = "0.689" }
    @{ Cell = "E23"; Value = "  -1.48%  " }
    @{ Cell = "D24"; Value = "83.64" }
    @{ Cell = "E25"; Value = "  -3.49%  " }
    @{ Cell = "D26"; Value = "2.24" }
    @{ Cell = "E26"; Value = "  -1.68%  " }
    @{ Cell = "D27"; Value = "10.09" }
    @{ Cell = "E27"; Value = "  -3.95%  " }
    @{ Cell = "D29"; Value = "7.89" }
    @{ Cell = "E29"; Value = "  +0.23%  " }
    @{ Cell = "E30"; Value = "  -1.90%  " }
    @{ Cell = "E31"; Value = "  -2.02%  " }
    @{ Cell = "D32"; Value = "27.99" }
    @{ Cell = "E32"; Value = "  -0.57%  " }
    @{ Cell = "B33"; Value = "PEPE" }
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" }
    @{ Cell = "D33"; Value = "0.0₃0942" }
    @{ Cell = "E33"; Value = "  +1.65%  " }
    @{ Cell = "B34"; Value = "Hedera" }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar" }
    @{ Cell = "D34"; Value = "0.111" }
    @{ Cell = "E34"; Value = "  -2.13%  " }
    @{ Cell = "E35"; Value = "  +0.04%  " }
    @{ Cell = "D36"; Value = "48.37" }
    @{ Cell = "E36"; Value = "  +3.11%  " }
    @{ Cell = "E37"; Value = "  -3.42%  " }
    @{ Cell = "D38"; Value = "0.942" }
    @{ Cell = "E38"; Value = "  -3.01%  " }
    @{ Cell = "E39"; Value = "  +1.64%  " }
    @{ Cell = "D40"; Value = "49.06" }
    @{ Cell = "E40"; Value = "  -1.94%  " }
    @{ Cell = "E41"; Value = "  -0.72%  " }
    @{ Cell = "E42"; Value = "  -0.34%  " }
    @{ Cell = "D43"; Value = "8.32" }
    @{ Cell = "E43"; Value = "  -1.41%  " }
    @{ Cell = "E44"; Value = "  +5.25%  " }
    @{ Cell = "D45"; Value = "2.798.28" }
    @{ Cell = "D46"; Value = "372.95" }
    @{ Cell = "E46"; Value = "  -2.94%  " }
    @{ Cell = "B47"; Value = "VeChain" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "D47"; Value = "0.0344" }
    @{ Cell = "E47"; Value = "  -1.56%  " }
    @{ Cell = "B48"; Value = "Monero" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" }
    @{ Cell = "D48"; Value = "135.73" }
    @{ Cell = "E48"; Value = "  +0.52%  " }
    @{ Cell = "D50"; Value = "25.73" }
    @{ Cell = "E51"; Value = "  +2.08%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $value = $u.Value

    # Price strings such as "576.41" or "0.152" round-trip as plain numbers
    # through COM's Range.Value setter, which would silently flip the cell
    # from text to a numeric type. The workbook stores every Price/Volume
    # value as text (others, like "66.960.98" or "  -0.13%  ", already fail
    # numeric parsing and are left alone), so a leading apostrophe forces
    # text for the ones that would otherwise be auto-converted; the
    # resulting quote-prefix style is then reset so formatting stays
    # identical to the original cell.
    if ($value -match '^\d+(\.\d+)?$') {
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
